# Add " (Mikayla)" right after "ARIMA model", with "(Mikayla)" rendered in
# bold — mirroring the existing "<model> (<Name>)" author-tag pattern
# already used for the other bullet points in this list (e.g.
# "Neural Net (Alexis)", "TSLM model (Cora)").

$d = $word.ActiveDocument

# Locate the "ARIMA model" run.
$rng = $d.Content
$found = $rng.Find.Execute("ARIMA model", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find 'ARIMA model' text"
}

# Collapse to the point right after "ARIMA model" and type the new text
# (space + author tag) there, inheriting the surrounding (non-bold)
# character formatting.
$rng.Collapse(0)
$rng.InsertAfter(" (Mikayla)")

# $rng now spans the newly inserted " (Mikayla)" text. Re-address just the
# "(Mikayla)" portion (the last 9 characters) and make it bold, leaving the
# leading space un-bolded.
$tagStart = $rng.End - 9
$tag = $d.Range($tagStart, $rng.End)
$tag.Font.Bold = 1
